$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = "미분방정식을 이용한 현상 모델링"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/01/modeling_with_differential_equation.html"

# Row 6
$ws.Range("D6").Value = "[R Markdown] Markdown, DT datatable, dygraph 기본 옵션"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/R-Markdown-Markdown-DT-datatable-dygraph-%EA%B8%B0%EB%B3%B8-%EC%98%B5%EC%85%98"

# Row 9
$ws.Range("D9").Value = "[공지] 대학원 입시 설명회 – 4월 28일 저녁 7시 + 후기"

# Row 23
$ws.Range("D23").Value = "2020년 가을에 UMASS에서 개설된 Advanced NLP 강의입니다.`n슬라이드/동영상 모두 제공됩니다.`n강의 제목처럼 기본 NLP내용 외"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2798"

# Row 51
$ws.Range("D51").Value = "[python] 딕셔너리 객체의 setdefault() 메소드 사용법"
$ws.Range("E51").Value = "https://bskyvision.com/1186"
